$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 106712.14
$ws.Range("I40").Value = 601080
$ws.Range("J40").Value = 3718.8333
$ws.Range("K40").Value = 601080
$ws.Range("L40").Value = 3718.8333
$ws.Range("M40").Value = -600905
$ws.Range("N40").Value = -4068.8333

$ws.Range("H62").Value = 9387.777
$ws.Range("I62").Value = 8000
$ws.Range("J62").Value = 10081.667
$ws.Range("K62").Value = 8000
$ws.Range("L62").Value = 10081.667
$ws.Range("M62").Value = -7376
$ws.Range("N62").Value = -11329.667

$ws.Range("H65").Value = 9387.777
$ws.Range("I65").Value = 8000
$ws.Range("J65").Value = 10081.667
$ws.Range("K65").Value = 40000
$ws.Range("L65").Value = 50408.335
$ws.Range("M65").Value = -36880
$ws.Range("N65").Value = -56648.335

$ws.Range("H113").Value = 2469
$ws.Range("I113").Value = 2082.889
$ws.Range("K113").Value = 2082.889
$ws.Range("M113").Value = 1171.111

$ws.Range("H132").Value = 4464.6045
$ws.Range("I132").Value = 4451.881
$ws.Range("K132").Value = 13355.643
$ws.Range("M132").Value = -10825.643

$ws.Range("H138").Value = 1967.0233
$ws.Range("I138").Value = 1160.7059
$ws.Range("J138").Value = 2494.2307
$ws.Range("K138").Value = 3482.1177
$ws.Range("L138").Value = 7482.6921
$ws.Range("M138").Value = 1657.8823
$ws.Range("N138").Value = -17762.6921

$ws.Range("H141").Value = 1124.25
$ws.Range("I141").Value = 1073
$ws.Range("K141").Value = 3219
$ws.Range("M141").Value = 1961

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15955927
$ws.Range("I32").Value = 15838599
$ws.Range("K32").Value = 15838599
$ws.Range("M32").Value = -15838312

$ws.Range("H60").Value = 29525
$ws.Range("I60").Value = 29525
$ws.Range("K60").Value = 29525
$ws.Range("M60").Value = -28792

$ws.Range("H61").Value = 1944.8868
$ws.Range("I61").Value = 1799.78
$ws.Range("K61").Value = 1799.78
$ws.Range("M61").Value = -1587.78

$ws.Range("H63").Value = 3479
$ws.Range("I63").Value = 2493.5
$ws.Range("K63").Value = 2493.5
$ws.Range("M63").Value = -1807.5

$ws.Range("H66").Value = 3479
$ws.Range("I66").Value = 2493.5
$ws.Range("K66").Value = 12467.5
$ws.Range("M66").Value = -9035.5

$ws.Range("H74").Value = 2633.2188
$ws.Range("I74").Value = 2663.1072
$ws.Range("J74").Value = 2424
$ws.Range("K74").Value = 2663.1072
$ws.Range("L74").Value = 2424
$ws.Range("M74").Value = -1789.1072
$ws.Range("N74").Value = -4172

$ws.Range("H77").Value = 2633.2188
$ws.Range("I77").Value = 2663.1072
$ws.Range("J77").Value = 2424
$ws.Range("K77").Value = 13315.536
$ws.Range("L77").Value = 12120
$ws.Range("M77").Value = -8947.536
$ws.Range("N77").Value = -20856

$ws.Range("H110").Value = 1887.4166
$ws.Range("I110").Value = 1575
$ws.Range("K110").Value = 1575
$ws.Range("M110").Value = 470

$ws.Range("H122").Value = 5738
$ws.Range("I122").Value = 3722.111
$ws.Range("J122").Value = 7249.9165
$ws.Range("K122").Value = 11166.333
$ws.Range("L122").Value = 21749.7495
$ws.Range("M122").Value = -8716.332999999999
$ws.Range("N122").Value = -26649.7495

$ws.Range("H132").Value = 2434.814
$ws.Range("I132").Value = 1737.4839
$ws.Range("K132").Value = 5212.4517
$ws.Range("M132").Value = -2682.4517

$ws.Range("H136").Value = 1944.8868
$ws.Range("I136").Value = 1799.78
$ws.Range("K136").Value = 5399.34
$ws.Range("M136").Value = -2849.34

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 96.333336
$ws.Range("J22").Value = 99
$ws.Range("L22").Value = 99
$ws.Range("N22").Value = -445

$ws.Range("H107").Value = 2562.9565
$ws.Range("I107").Value = 2508.8235
$ws.Range("K107").Value = 2508.8235
$ws.Range("M107").Value = -588.8235

$ws.Range("H112").Value = 75329.664
$ws.Range("J112").Value = 75329.664
$ws.Range("L112").Value = 75329.664
$ws.Range("N112").Value = -78283.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2677.862
$ws.Range("I58").Value = 2676.2964
$ws.Range("J58").Value = 2699
$ws.Range("K58").Value = 2676.2964
$ws.Range("L58").Value = 2699
$ws.Range("M58").Value = -2473.2964
$ws.Range("N58").Value = -3105

$ws.Range("H68").Value = 51348.125
$ws.Range("J68").Value = 62957
$ws.Range("L68").Value = 62957
$ws.Range("N68").Value = -64455

$ws.Range("H71").Value = 51348.125
$ws.Range("J71").Value = 62957
$ws.Range("L71").Value = 188871
$ws.Range("N71").Value = -196359

$ws.Range("H122").Value = 3849139
$ws.Range("I122").Value = 5003243.5
$ws.Range("J122").Value = 2124.6667
$ws.Range("K122").Value = 15009730.5
$ws.Range("L122").Value = 6374.000100000001
$ws.Range("M122").Value = -15007280.5
$ws.Range("N122").Value = -11274.0001

$ws.Range("H132").Value = 3634.5386
$ws.Range("I132").Value = 3481.2046
$ws.Range("K132").Value = 10443.6138
$ws.Range("M132").Value = -7913.613799999999

$ws.Range("H134").Value = 2898.5386
$ws.Range("I134").Value = 2653.5
$ws.Range("J134").Value = 4246.25
$ws.Range("K134").Value = 7960.5
$ws.Range("L134").Value = 12738.75
$ws.Range("M134").Value = -5425.5
$ws.Range("N134").Value = -17808.75

$ws.Range("H136").Value = 2677.862
$ws.Range("I136").Value = 2676.2964
$ws.Range("J136").Value = 2699
$ws.Range("K136").Value = 8028.889200000001
$ws.Range("L136").Value = 8097
$ws.Range("M136").Value = -5478.889200000001
$ws.Range("N136").Value = -13197

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2330.6667
$ws.Range("I64").Value = 2330.6667
$ws.Range("K64").Value = 6992.000100000001
$ws.Range("M64").Value = -6722.000100000001

$ws.Range("H67").Value = 2330.6667
$ws.Range("I67").Value = 2330.6667
$ws.Range("K67").Value = 6992.000100000001
$ws.Range("M67").Value = -6056.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 10000
$ws.Range("I59").Value = 10000
$ws.Range("K59").Value = 10000
$ws.Range("M59").Value = -9417

$ws.Range("H113").Value = 1996.5
$ws.Range("I113").Value = 1996.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1996.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 173.5
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 2281.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2281.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6844.5
$ws.Range("N122").Value = -11744.5
$ws.Range("M122").ClearContents()

$ws.Range("H125").Value = 85055.60000000001
$ws.Range("J125").Value = 85055.60000000001
$ws.Range("L125").Value = 85055.60000000001
$ws.Range("N125").Value = -89975.60000000001

$ws.Range("H132").Value = 4207.073
$ws.Range("I132").Value = 3961.2368
$ws.Range("K132").Value = 11883.7104
$ws.Range("M132").Value = -9353.7104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 41675370
$ws.Range("I40").Value = 55562156
$ws.Range("K40").Value = 55562156
$ws.Range("M40").Value = -55562020

$ws.Range("H61").Value = 4642.5
$ws.Range("I61").Value = 3034.25
$ws.Range("J61").Value = 20725
$ws.Range("K61").Value = 3034.25
$ws.Range("L61").Value = 20725
$ws.Range("M61").Value = -2832.25
$ws.Range("N61").Value = -21129

$ws.Range("H75").Value = 109994.5
$ws.Range("J75").Value = 109994.5
$ws.Range("L75").Value = 109994.5
$ws.Range("N75").Value = -111866.5

$ws.Range("H78").Value = 109994.5
$ws.Range("J78").Value = 109994.5
$ws.Range("L78").Value = 329983.5
$ws.Range("N78").Value = -339343.5

$ws.Range("H113").Value = 4642.5
$ws.Range("I113").Value = 3034.25
$ws.Range("J113").Value = 20725
$ws.Range("K113").Value = 3034.25
$ws.Range("L113").Value = 20725
$ws.Range("M113").Value = -864.25
$ws.Range("N113").Value = -25065

$ws.Range("H122").Value = 57867.125
$ws.Range("I122").Value = 59278.145
$ws.Range("J122").Value = 47990
$ws.Range("K122").Value = 177834.435
$ws.Range("L122").Value = 143970
$ws.Range("M122").Value = -175384.435
$ws.Range("N122").Value = -148870

$ws.Range("H132").Value = 27642.38
$ws.Range("I132").Value = 30595
$ws.Range("K132").Value = 91785
$ws.Range("M132").Value = -89255

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 12516319
$ws.Range("I2").Value = 12516319
$ws.Range("K2").Value = 12516319
$ws.Range("M2").Value = -12516207

$ws.Range("H122").Value = 200008480
$ws.Range("I122").Value = 500003500
$ws.Range("J122").Value = 11792.333
$ws.Range("K122").Value = 1500010500
$ws.Range("L122").Value = 35376.999
$ws.Range("M122").Value = -1500008050
$ws.Range("N122").Value = -40276.999

$ws.Range("H132").Value = 2546.7222
$ws.Range("I132").Value = 2056.9333
$ws.Range("K132").Value = 6170.7999
$ws.Range("M132").Value = -3640.7999

$ws.Range("H136").Value = 39352.89
$ws.Range("I136").Value = 1744.6666
$ws.Range("K136").Value = 5233.9998
$ws.Range("M136").Value = -2683.9998
